# Apply updated Betfair Back/Lay odds values for Jogos_do_Dia workbook
# Targets the single worksheet in the workbook ("Sheet1")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 2.32
$ws.Range("G2").Value = 2.62
$ws.Range("H2").Value = 3.3
$ws.Range("K2").Value = 3.8
$ws.Range("N2").Value = 2.64
$ws.Range("O2").Value = 1.42
$ws.Range("P2").Value = 1.65
$ws.Range("Q2").Value = 2.12
$ws.Range("R2").Value = 1.24
$ws.Range("U2").Value = 1.89
$ws.Range("G3").Value = 1.25
$ws.Range("P3").Value = 2.16
$ws.Range("Q3").Value = 1.7
$ws.Range("F4").Value = 1.72
$ws.Range("G4").Value = 1.89
$ws.Range("H4").Value = 4.3
$ws.Range("I4").Value = 5.2
$ws.Range("J4").Value = 3.95
$ws.Range("K4").Value = 4.7
$ws.Range("P4").Value = 2.3
$ws.Range("Q4").Value = 1.58
$ws.Range("O5").Value = 1.54
$ws.Range("P5").Value = 1.57
$ws.Range("X5").Value = 8.6
$ws.Range("AN5").Value = 27
$ws.Range("F8").Value = 2.02
$ws.Range("G8").Value = 2.28
$ws.Range("H8").Value = 3.75
$ws.Range("I8").Value = 5.1
$ws.Range("J8").Value = 2.96
$ws.Range("K8").Value = 3.75
$ws.Range("N8").Value = 3.15
$ws.Range("S8").Value = 3.7
$ws.Range("T8").Value = 1.83
$ws.Range("U8").Value = 1.96
$ws.Range("V8").Value = 1.24
$ws.Range("W8").Value = 1.78
$ws.Range("Y8").Value = 16.5
$ws.Range("Z8").Value = 980
$ws.Range("AA8").Value = 110
$ws.Range("AC8").Value = 9.2
$ws.Range("AD8").Value = 980
$ws.Range("AF8").Value = 15.5
$ws.Range("AG8").Value = 13
$ws.Range("AH8").Value = 980
$ws.Range("AI8").Value = 80
$ws.Range("AJ8").Value = 980
$ws.Range("AK8").Value = 980
$ws.Range("AM8").Value = 140
$ws.Range("AN8").Value = 23
$ws.Range("AO8").Value = 80
$ws.Range("F9").Value = 2.96
$ws.Range("G9").Value = 3.35
$ws.Range("H9").Value = 2.7
$ws.Range("I9").Value = 2.94
$ws.Range("K9").Value = 3.25
$ws.Range("P9").Value = 1.55
$ws.Range("J10").Value = 6.4
$ws.Range("U10").Value = 1.68
$ws.Range("AF10").Value = 7.2
$ws.Range("P11").Value = 1.84
$ws.Range("Q11").Value = 2.12
$ws.Range("S11").Value = 3.8
$ws.Range("U11").Value = 2.1
$ws.Range("X11").Value = 13
$ws.Range("Y11").Value = 11
$ws.Range("Z11").Value = 18.5
$ws.Range("AA11").Value = 46
$ws.Range("AB11").Value = 11
$ws.Range("AD11").Value = 13
$ws.Range("AE11").Value = 34
$ws.Range("AF11").Value = 18.5
$ws.Range("AG11").Value = 13
$ws.Range("AH11").Value = 18
$ws.Range("AK11").Value = 34
$ws.Range("AM11").Value = 130
$ws.Range("AN11").Value = 30
$ws.Range("AO11").Value = 32
